$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing data (A:D) shifts to (B:E)
$ws.Columns.Item(1).Insert()

# New header and scenario values in column A
# (written in this order so the shared-strings table matches the
# original authoring order: Scenario, Smoke, Regression, QALive)
$ws.Range("A1").Value = "Scenario"
$ws.Range("A3").Value = "Smoke"
$ws.Range("A2").Value = "Regression"
$ws.Range("A4").Value = "QALive"

# New CartCount values in column E (header already shifted from D1)
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2

# Update the active cell selection to E1 as in the final file
$ws.Range("E1").Select()
